{"js": "// Replace the date and each \"NNN\u00f7N=\" division problem with its updated value.\n// Each old string is unique in the document, so a simple search + full-text\n// replace on the single match found for each is unambiguous and safe.\nconst replacements = [\n  [\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"],\n  [\"499\u00f75=\", \"337\u00f73=\"],\n  [\"585\u00f72=\", \"806\u00f76=\"],\n  [\"340\u00f74=\", \"573\u00f77=\"],\n  [\"944\u00f76=\", \"782\u00f75=\"],\n  [\"429\u00f78=\", \"529\u00f73=\"],\n  [\"915\u00f78=\", \"659\u00f74=\"],\n  [\"807\u00f79=\", \"681\u00f78=\"],\n  [\"462\u00f79=\", \"842\u00f76=\"],\n  [\"871\u00f74=\", \"971\u00f72=\"],\n  [\"729\u00f73=\", \"335\u00f78=\"],\n  [\"299\u00f73=\", \"289\u00f76=\"],\n  [\"125\u00f73=\", \"788\u00f78=\"],\n  [\"767\u00f77=\", \"828\u00f77=\"],\n  [\"778\u00f74=\", \"221\u00f78=\"],\n  [\"492\u00f74=\", \"406\u00f74=\"],\n  [\"695\u00f77=\", \"954\u00f77=\"],\n  [\"538\u00f72=\", \"285\u00f75=\"],\n  [\"571\u00f73=\", \"568\u00f73=\"],\n  [\"188\u00f77=\", \"502\u00f74=\"],\n  [\"793\u00f79=\", \"870\u00f77=\"],\n  [\"264\u00f73=\", \"623\u00f79=\"],\n  [\"239\u00f77=\", \"405\u00f79=\"],\n  [\"417\u00f74=\", \"234\u00f78=\"],\n  [\"290\u00f75=\", \"834\u00f72=\"],\n  [\"433\u00f76=\", \"289\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each \"NNN\u00f7N=\" division problem with its updated\n# value. Each old string is unique in the document, so Find/Replace with\n# wdReplaceAll (but only ever matching a single hit) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-12 Monday\", \"2024-02-13 Tuesday\"),\n    @(\"499\u00f75=\", \"337\u00f73=\"),\n    @(\"585\u00f72=\", \"806\u00f76=\"),\n    @(\"340\u00f74=\", \"573\u00f77=\"),\n    @(\"944\u00f76=\", \"782\u00f75=\"),\n    @(\"429\u00f78=\", \"529\u00f73=\"),\n    @(\"915\u00f78=\", \"659\u00f74=\"),\n    @(\"807\u00f79=\", \"681\u00f78=\"),\n    @(\"462\u00f79=\", \"842\u00f76=\"),\n    @(\"871\u00f74=\", \"971\u00f72=\"),\n    @(\"729\u00f73=\", \"335\u00f78=\"),\n    @(\"299\u00f73=\", \"289\u00f76=\"),\n    @(\"125\u00f73=\", \"788\u00f78=\"),\n    @(\"767\u00f77=\", \"828\u00f77=\"),\n    @(\"778\u00f74=\", \"221\u00f78=\"),\n    @(\"492\u00f74=\", \"406\u00f74=\"),\n    @(\"695\u00f77=\", \"954\u00f77=\"),\n    @(\"538\u00f72=\", \"285\u00f75=\"),\n    @(\"571\u00f73=\", \"568\u00f73=\"),\n    @(\"188\u00f77=\", \"502\u00f74=\"),\n    @(\"793\u00f79=\", \"870\u00f77=\"),\n    @(\"264\u00f73=\", \"623\u00f79=\"),\n    @(\"239\u00f77=\", \"405\u00f79=\"),\n    @(\"417\u00f74=\", \"234\u00f78=\"),\n    @(\"290\u00f75=\", \"834\u00f72=\"),\n    @(\"433\u00f76=\", \"289\u00f76=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute([ref]$old, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$new, [ref]2)\n}\n"}
